$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1001
$ws.Range("I18").Value = 1001
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1001
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -717
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H51").Value = 3404.6667
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 3580.25
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 3580.25
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -4548.25
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value = 0
$ws.Range("H100").Value = 1852.0968
$ws.Range("I100").Value = 1108.0435
$ws.Range("J100").Value = 3991.25
$ws.Range("K100").Value = 1108.0435
$ws.Range("L100").Value = 3991.25
$ws.Range("M100").Value = -567.0435
$ws.Range("N100").Value = -5073.25
$ws.Range("H133").Value = 59800
$ws.Range("J133").Value = 59800
$ws.Range("L133").Value = 59800
$ws.Range("N133").Value = -69920
$ws.Range("H138").Value = 1893.6923
$ws.Range("I138").Value = 1409.375
$ws.Range("J138").Value = 3508.0833
$ws.Range("K138").Value = 4228.125
$ws.Range("L138").Value = 10524.2499
$ws.Range("M138").Value = 911.875
$ws.Range("N138").Value = -20804.2499
$ws.Range("H139").Value = 70565
$ws.Range("J139").Value = 70565
$ws.Range("L139").Value = 70565
$ws.Range("N139").Value = -80845
$ws.Range("H140").Value = 80936
$ws.Range("J140").Value = 80936
$ws.Range("L140").Value = 80936
$ws.Range("N140").Value = -91296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 6679.6665
$ws.Range("I37").Value = 2400
$ws.Range("J37").Value = 8106.222
$ws.Range("K37").Value = 2400
$ws.Range("L37").Value = 8106.222
$ws.Range("M37").Value = -2127
$ws.Range("N37").Value = -8652.222
$ws.Range("H61").Value = 2673.054
$ws.Range("I61").Value = 3109.7917
$ws.Range("J61").Value = 1866.7693
$ws.Range("K61").Value = 3109.7917
$ws.Range("L61").Value = 1866.7693
$ws.Range("M61").Value = -2897.7917
$ws.Range("N61").Value = -2290.7693
$ws.Range("H74").Value = 1310.1482
$ws.Range("I74").Value = 901.8946999999999
$ws.Range("J74").Value = 2279.75
$ws.Range("K74").Value = 901.8946999999999
$ws.Range("L74").Value = 2279.75
$ws.Range("M74").Value = -27.89469999999994
$ws.Range("N74").Value = -4027.75
$ws.Range("H77").Value = 1310.1482
$ws.Range("I77").Value = 901.8946999999999
$ws.Range("J77").Value = 2279.75
$ws.Range("K77").Value = 4509.4735
$ws.Range("L77").Value = 11398.75
$ws.Range("M77").Value = -141.4735000000001
$ws.Range("N77").Value = -20134.75
$ws.Range("H88").Value = 4881.7
$ws.Range("I88").Value = 2229.5
$ws.Range("J88").Value = 6649.8335
$ws.Range("K88").Value = 2229.5
$ws.Range("L88").Value = 6649.8335
$ws.Range("M88").Value = -1823.5
$ws.Range("N88").Value = -7461.8335
$ws.Range("H91").Value = 4881.7
$ws.Range("I91").Value = 2229.5
$ws.Range("J91").Value = 6649.8335
$ws.Range("K91").Value = 2229.5
$ws.Range("L91").Value = 6649.8335
$ws.Range("M91").Value = -825.5
$ws.Range("N91").Value = -9457.833500000001
$ws.Range("H97").Value = 1452.2565
$ws.Range("I97").Value = 1181.6333
$ws.Range("J97").Value = 2354.3333
$ws.Range("K97").Value = 1181.6333
$ws.Range("L97").Value = 2354.3333
$ws.Range("M97").Value = -685.6333
$ws.Range("N97").Value = -3346.3333
$ws.Range("H102").Value = 4103.793
$ws.Range("I102").Value = 4242.7856
$ws.Range("J102").Value = 3974.0667
$ws.Range("K102").Value = 4242.7856
$ws.Range("L102").Value = 3974.0667
$ws.Range("M102").Value = -2620.7856
$ws.Range("N102").Value = -7218.066699999999
$ws.Range("H132").Value = 2397.2173
$ws.Range("I132").Value = 1901.3334
$ws.Range("J132").Value = 4182.4
$ws.Range("K132").Value = 5704.0002
$ws.Range("L132").Value = 12547.2
$ws.Range("M132").Value = -3174.0002
$ws.Range("N132").Value = -17607.2
$ws.Range("H136").Value = 2673.054
$ws.Range("I136").Value = 3109.7917
$ws.Range("J136").Value = 1866.7693
$ws.Range("K136").Value = 9329.375100000001
$ws.Range("L136").Value = 5600.3079
$ws.Range("M136").Value = -6779.375100000001
$ws.Range("N136").Value = -10700.3079
$ws.Range("H138").Value = 50368.57
$ws.Range("J138").Value = 50368.57
$ws.Range("L138").Value = 50368.57
$ws.Range("N138").Value = -60648.57
$ws.Range("H141").Value = 63411.11
$ws.Range("J141").Value = 63411.11
$ws.Range("L141").Value = 63411.11
$ws.Range("N141").Value = -73771.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1667275.5
$ws.Range("I7").Value = 663.25
$ws.Range("J7").Value = 5000500
$ws.Range("K7").Value = 663.25
$ws.Range("L7").Value = 5000500
$ws.Range("M7").Value = -550.25
$ws.Range("N7").Value = -5000726
$ws.Range("H86").Value = 3563.0833
$ws.Range("I86").Value = 3655.077
$ws.Range("J86").Value = 3454.3635
$ws.Range("K86").Value = 3655.077
$ws.Range("L86").Value = 3454.3635
$ws.Range("M86").Value = -2532.077
$ws.Range("N86").Value = -5700.363499999999
$ws.Range("H89").Value = 3563.0833
$ws.Range("I89").Value = 3655.077
$ws.Range("J89").Value = 3454.3635
$ws.Range("K89").Value = 18275.385
$ws.Range("L89").Value = 17271.8175
$ws.Range("M89").Value = -12659.385
$ws.Range("N89").Value = -28503.8175
$ws.Range("H94").Value = 848.1177
$ws.Range("I94").Value = 699.2593000000001
$ws.Range("J94").Value = 1422.2858
$ws.Range("K94").Value = 699.2593000000001
$ws.Range("L94").Value = 1422.2858
$ws.Range("M94").Value = -248.2593000000001
$ws.Range("N94").Value = -2324.2858
$ws.Range("H99").Value = 1912.375
$ws.Range("I99").Value = 1518.7
$ws.Range("J99").Value = 2568.5
$ws.Range("K99").Value = 1518.7
$ws.Range("L99").Value = 2568.5
$ws.Range("M99").Value = -20.70000000000005
$ws.Range("N99").Value = -5564.5
$ws.Range("H134").Value = 5157.914
$ws.Range("I134").Value = 650.1923
$ws.Range("J134").Value = 18180.223
$ws.Range("K134").Value = 1950.5769
$ws.Range("L134").Value = 54540.66900000001
$ws.Range("M134").Value = 584.4231
$ws.Range("N134").Value = -59610.66900000001
$ws.Range("H140").Value = 89633.336
$ws.Range("J140").Value = 89633.336
$ws.Range("L140").Value = 89633.336
$ws.Range("N140").Value = -99993.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 3733.25
$ws.Range("J8").Value = 4000
$ws.Range("L8").Value = 4000
$ws.Range("N8").Value = -4280
$ws.Range("H56").Value = 7100
$ws.Range("J56").Value = 9200
$ws.Range("L56").Value = 9200
$ws.Range("N56").Value = -10890
$ws.Range("H60").Value = 8300.75
$ws.Range("J60").Value = 8300.75
$ws.Range("L60").Value = 8300.75
$ws.Range("N60").Value = -9322.75
$ws.Range("H132").Value = 2300.6453
$ws.Range("I132").Value = 1882.1818
$ws.Range("J132").Value = 3323.5557
$ws.Range("K132").Value = 5646.5454
$ws.Range("L132").Value = 9970.667099999999
$ws.Range("M132").Value = -3116.5454
$ws.Range("N132").Value = -15030.6671
$ws.Range("H134").Value = 28572716
$ws.Range("I134").Value = 40001404
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 120004212
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -120001677
$ws.Range("N134").Value = -8070
$ws.Range("H140").Value = 87390
$ws.Range("J140").Value = 87390
$ws.Range("L140").Value = 87390
$ws.Range("N140").Value = -97750

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 147
$ws.Range("J50").Value = 300
$ws.Range("L50").Value = 900
$ws.Range("N50").Value = -1862
$ws.Range("H53").Value = 147
$ws.Range("J53").Value = 300
$ws.Range("L53").Value = 900
$ws.Range("N53").Value = -1862
$ws.Range("H69").Value = 1527.2727
$ws.Range("J69").Value = 1800
$ws.Range("L69").Value = 5400
$ws.Range("N69").Value = -7022
$ws.Range("H72").Value = 1527.2727
$ws.Range("J72").Value = 1800
$ws.Range("L72").Value = 16200
$ws.Range("N72").Value = -24312
$ws.Range("H98").Value = 1007.7143
$ws.Range("I98").Value = 723.3333
$ws.Range("J98").Value = 1221
$ws.Range("K98").Value = 2169.9999
$ws.Range("L98").Value = 3663
$ws.Range("M98").Value = -671.9998999999998
$ws.Range("N98").Value = -6659
$ws.Range("H107").Value = 620.2222
$ws.Range("J107").Value = 654
$ws.Range("L107").Value = 1962
$ws.Range("N107").Value = -5802
$ws.Range("H132").Value = 1570.9375
$ws.Range("I132").Value = 1264.2307
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 11378.0763
$ws.Range("L132").Value = 26100
$ws.Range("M132").Value = -8848.076300000001
$ws.Range("N132").Value = -31160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 142859710
$ws.Range("J33").Value = 153848690
$ws.Range("L33").Value = 153848690
$ws.Range("N33").Value = -153849194
$ws.Range("H53").Value = 53000
$ws.Range("I53").Value = 8000
$ws.Range("K53").Value = 8000
$ws.Range("M53").Value = -7369
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").ClearContents()
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = 0
$ws.Range("H97").Value = 2682
$ws.Range("I97").Value = 950
$ws.Range("J97").Value = 4414
$ws.Range("K97").Value = 950
$ws.Range("L97").Value = 4414
$ws.Range("M97").Value = -454
$ws.Range("N97").Value = -5406
$ws.Range("H133").Value = 59520
$ws.Range("J133").Value = 59520
$ws.Range("L133").Value = 59520
$ws.Range("N133").Value = -69640
$ws.Range("H140").Value = 41769.75
$ws.Range("J140").Value = 41769.75
$ws.Range("L140").Value = 41769.75
$ws.Range("N140").Value = -52129.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 583.3333
$ws.Range("I9").Value = 375
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 375
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = -151
$ws.Range("N9").Value = -1448
$ws.Range("H55").Value = 682.5484
$ws.Range("I55").Value = 823.1667
$ws.Range("J55").Value = 487.84616
$ws.Range("K55").Value = 823.1667
$ws.Range("L55").Value = 487.84616
$ws.Range("M55").Value = -650.1667
$ws.Range("N55").Value = -833.8461600000001
$ws.Range("H82").Value = 2485.7144
$ws.Range("I82").Value = 3466.6667
$ws.Range("J82").Value = 1750
$ws.Range("K82").Value = 3466.6667
$ws.Range("L82").Value = 1750
$ws.Range("M82").Value = -3105.6667
$ws.Range("N82").Value = -2472
$ws.Range("H85").Value = 2485.7144
$ws.Range("I85").Value = 3466.6667
$ws.Range("J85").Value = 1750
$ws.Range("K85").Value = 3466.6667
$ws.Range("L85").Value = 1750
$ws.Range("M85").Value = -2218.6667
$ws.Range("N85").Value = -4246
$ws.Range("H100").Value = 3543.1052
$ws.Range("I100").Value = 3283.1428
$ws.Range("J100").Value = 3694.75
$ws.Range("K100").Value = 3283.1428
$ws.Range("L100").Value = 3694.75
$ws.Range("M100").Value = -2742.1428
$ws.Range("N100").Value = -4776.75
$ws.Range("H132").Value = 2640.8364
$ws.Range("I132").Value = 2741.2285
$ws.Range("J132").Value = 2465.15
$ws.Range("K132").Value = 8223.6855
$ws.Range("L132").Value = 7395.450000000001
$ws.Range("M132").Value = -5693.6855
$ws.Range("N132").Value = -12455.45
$ws.Range("H133").Value = 92335.71000000001
$ws.Range("J133").Value = 92335.71000000001
$ws.Range("L133").Value = 92335.71000000001
$ws.Range("N133").Value = -97395.71000000001
$ws.Range("H138").Value = 63322.43
$ws.Range("J138").Value = 63322.43
$ws.Range("L138").Value = 63322.43
$ws.Range("N138").Value = -73602.42999999999
$ws.Range("H139").Value = 79800
$ws.Range("J139").Value = 79800
$ws.Range("L139").Value = 79800
$ws.Range("N139").Value = -90080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 640
$ws.Range("I4").Value = 400
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -287
$ws.Range("N4").Value = -1226
$ws.Range("H96").Value = 3515.7896
$ws.Range("I96").Value = 2225
$ws.Range("J96").Value = 5728.5713
$ws.Range("K96").Value = 2225
$ws.Range("L96").Value = 5728.5713
$ws.Range("M96").Value = -852
$ws.Range("N96").Value = -8474.5713
$ws.Range("H132").Value = 100002380
$ws.Range("I132").Value = 125001980
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 375005940
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -375003410
$ws.Range("N132").Value = -17057.9999
$ws.Range("H139").Value = 54057.145
$ws.Range("J139").Value = 54057.145
$ws.Range("L139").Value = 54057.145
$ws.Range("N139").Value = -64337.145
